$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header column in H1, matching the style/formatting of the
# existing header cells (e.g. G1: bold, bordered, centered).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add the new data values in H2 and H3
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
